$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: add numeric values in C1:F1
$ws.Range("C1").Value = 2
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 4
$ws.Range("F1").Value = 5

# Row 2: add numeric value in F2
$ws.Range("F2").Value = 2321

# Row 3: new row, all numeric values
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 6

# Move selection to G3 to match the final saved state
$ws.Range("G3").Select()
